# Apply integer ("0") number format to the Growth column (G) for the
# existing data rows as well as three extra blank rows appended below the
# table (rows 20-22), mirroring the new "shiny interactive app" export that
# widened the working range and reformatted the numeric column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Extend formatting (and therefore the sheet's used range) down to row 22,
# applying the built-in "0" number format (numFmtId 1) to G2:G22.
$ws.Range("G2:G22").NumberFormat = "0"

# Move/restore the active selection to C9, matching the saved workbook view.
$ws.Range("C9").Select()
